# Change the table style applied to the table on slide 6 (the data table)
# from the custom theme-linked "Table_0" style to the built-in table
# style {573E8F96-C595-4247-B0C5-C318D0C3D19F}.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{573E8F96-C595-4247-B0C5-C318D0C3D19F}")
        }
    }
}
